$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix formatting of row 136 (last row of event 115) to match the
# "end of event group" style used by rows 10,19,...,127 (border + matching styles)
$null = $ws.Range("A127:G127").Copy()
$null = $ws.Range("A136:G136").PasteSpecial(-4122)

# --- Append two new 9-row event blocks (events 116 and 117), rows 137-154.
# Use existing rows 119:127 (event 114, a full normal+last-row template)
# as a formatting template, then overwrite the values.
$null = $ws.Range("A119:G127").Copy()
$null = $ws.Range("A137:G145").PasteSpecial(-4122)
$null = $ws.Range("A119:G127").Copy()
$null = $ws.Range("A146:G154").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Fill in the values for the new rows ---
$ws.Cells.Item(137,1).Value = 116
$ws.Cells.Item(137,2).NumberFormat = "@"
$ws.Cells.Item(137,2).Value = "08"
$ws.Cells.Item(137,3).Value = 1
$ws.Cells.Item(137,4).Value = "drawing"
$ws.Cells.Item(137,5).Value = 0
$ws.Cells.Item(137,6).Value = 10
$ws.Cells.Item(137,7).Value = 8

$ws.Cells.Item(138,1).Value = 116
$ws.Cells.Item(138,2).NumberFormat = "@"
$ws.Cells.Item(138,2).Value = "08"
$ws.Cells.Item(138,3).Value = 2
$ws.Cells.Item(138,4).Value = "whispering"
$ws.Cells.Item(138,5).Value = 8
$ws.Cells.Item(138,6).Value = 7
$ws.Cells.Item(138,7).Value = 7

$ws.Cells.Item(139,1).Value = 116
$ws.Cells.Item(139,2).NumberFormat = "@"
$ws.Cells.Item(139,2).Value = "08"
$ws.Cells.Item(139,3).Value = 3
$ws.Cells.Item(139,4).Value = "locking at phone"
$ws.Cells.Item(139,5).Value = 5
$ws.Cells.Item(139,6).Value = 6
$ws.Cells.Item(139,7).Value = 8

$ws.Cells.Item(140,1).Value = 116
$ws.Cells.Item(140,2).NumberFormat = "@"
$ws.Cells.Item(140,2).Value = "08"
$ws.Cells.Item(140,3).Value = 4
$ws.Cells.Item(140,4).Value = "heckling"
$ws.Cells.Item(140,5).Value = 10
$ws.Cells.Item(140,6).Value = 3
$ws.Cells.Item(140,7).Value = 3

$ws.Cells.Item(141,1).Value = 116
$ws.Cells.Item(141,2).NumberFormat = "@"
$ws.Cells.Item(141,2).Value = "08"
$ws.Cells.Item(141,3).Value = 5
$ws.Cells.Item(141,4).Value = "clicking pen"
$ws.Cells.Item(141,5).Value = 8
$ws.Cells.Item(141,6).Value = 7
$ws.Cells.Item(141,7).Value = 5

$ws.Cells.Item(142,1).Value = 116
$ws.Cells.Item(142,2).NumberFormat = "@"
$ws.Cells.Item(142,2).Value = "08"
$ws.Cells.Item(142,3).Value = 6
$ws.Cells.Item(142,4).Value = "snipping"
$ws.Cells.Item(142,5).Value = 0
$ws.Cells.Item(142,6).Value = 10
$ws.Cells.Item(142,7).Value = 5

$ws.Cells.Item(143,1).Value = 116
$ws.Cells.Item(143,2).NumberFormat = "@"
$ws.Cells.Item(143,2).Value = "08"
$ws.Cells.Item(143,3).Value = 7
$ws.Cells.Item(143,4).Value = "drumming"
$ws.Cells.Item(143,5).Value = 9
$ws.Cells.Item(143,6).Value = 10
$ws.Cells.Item(143,7).Value = 5

$ws.Cells.Item(144,1).Value = 116
$ws.Cells.Item(144,2).NumberFormat = "@"
$ws.Cells.Item(144,2).Value = "08"
$ws.Cells.Item(144,3).Value = 8
$ws.Cells.Item(144,4).Value = "head on table"
$ws.Cells.Item(144,5).Value = 8
$ws.Cells.Item(144,6).Value = 10
$ws.Cells.Item(144,7).Value = 7

$ws.Cells.Item(145,1).Value = 116
$ws.Cells.Item(145,2).NumberFormat = "@"
$ws.Cells.Item(145,2).Value = "08"
$ws.Cells.Item(145,3).Value = 9
$ws.Cells.Item(145,4).Value = "chatting"
$ws.Cells.Item(145,5).Value = 10
$ws.Cells.Item(145,6).Value = 8
$ws.Cells.Item(145,7).Value = 3

$ws.Cells.Item(146,1).Value = 117
$ws.Cells.Item(146,2).NumberFormat = "@"
$ws.Cells.Item(146,2).Value = "01"
$ws.Cells.Item(146,3).Value = 1
$ws.Cells.Item(146,4).Value = "whispering"
$ws.Cells.Item(146,5).Value = 3
$ws.Cells.Item(146,6).Value = 4
$ws.Cells.Item(146,7).Value = 10

$ws.Cells.Item(147,1).Value = 117
$ws.Cells.Item(147,2).NumberFormat = "@"
$ws.Cells.Item(147,2).Value = "01"
$ws.Cells.Item(147,3).Value = 2
$ws.Cells.Item(147,4).Value = "heckling"
$ws.Cells.Item(147,5).Value = 7
$ws.Cells.Item(147,6).Value = 3
$ws.Cells.Item(147,7).Value = 3

$ws.Cells.Item(148,1).Value = 117
$ws.Cells.Item(148,2).NumberFormat = "@"
$ws.Cells.Item(148,2).Value = "01"
$ws.Cells.Item(148,3).Value = 3
$ws.Cells.Item(148,4).Value = "drawing"
$ws.Cells.Item(148,5).Value = 1
$ws.Cells.Item(148,6).Value = 5
$ws.Cells.Item(148,7).Value = 9

$ws.Cells.Item(149,1).Value = 117
$ws.Cells.Item(149,2).NumberFormat = "@"
$ws.Cells.Item(149,2).Value = "01"
$ws.Cells.Item(149,3).Value = 4
$ws.Cells.Item(149,4).Value = "snipping"
$ws.Cells.Item(149,5).Value = 4
$ws.Cells.Item(149,6).Value = 6
$ws.Cells.Item(149,7).Value = 5

$ws.Cells.Item(150,1).Value = 117
$ws.Cells.Item(150,2).NumberFormat = "@"
$ws.Cells.Item(150,2).Value = "01"
$ws.Cells.Item(150,3).Value = 5
$ws.Cells.Item(150,4).Value = "locking at phone"
$ws.Cells.Item(150,5).Value = 4
$ws.Cells.Item(150,6).Value = 6
$ws.Cells.Item(150,7).Value = 9

$ws.Cells.Item(151,1).Value = 117
$ws.Cells.Item(151,2).NumberFormat = "@"
$ws.Cells.Item(151,2).Value = "01"
$ws.Cells.Item(151,3).Value = 6
$ws.Cells.Item(151,4).Value = "head on table"
$ws.Cells.Item(151,5).Value = 1
$ws.Cells.Item(151,6).Value = 4
$ws.Cells.Item(151,7).Value = 7

$ws.Cells.Item(152,1).Value = 117
$ws.Cells.Item(152,2).NumberFormat = "@"
$ws.Cells.Item(152,2).Value = "01"
$ws.Cells.Item(152,3).Value = 7
$ws.Cells.Item(152,4).Value = "clicking pen"
$ws.Cells.Item(152,5).Value = 5
$ws.Cells.Item(152,6).Value = 6
$ws.Cells.Item(152,7).Value = 6

$ws.Cells.Item(153,1).Value = 117
$ws.Cells.Item(153,2).NumberFormat = "@"
$ws.Cells.Item(153,2).Value = "01"
$ws.Cells.Item(153,3).Value = 8
$ws.Cells.Item(153,4).Value = "drumming"
$ws.Cells.Item(153,5).Value = 5
$ws.Cells.Item(153,6).Value = 5
$ws.Cells.Item(153,7).Value = 7

$ws.Cells.Item(154,1).Value = 117
$ws.Cells.Item(154,2).NumberFormat = "@"
$ws.Cells.Item(154,2).Value = "01"
$ws.Cells.Item(154,3).Value = 9
$ws.Cells.Item(154,4).Value = "chatting"
$ws.Cells.Item(154,5).Value = 7
$ws.Cells.Item(154,6).Value = 4
$ws.Cells.Item(154,7).Value = 7

# --- Update the view: active selection ---
$null = $ws.Range("F149").Select()

Write-Host "done"
